$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Every paragraph that row 10 onward needs in its final state, spelled out as
# single-quoted here-strings so accents, quotes, "&", "<", ">" and apostrophes
# all come through completely literally (no escaping surprises).
# ---------------------------------------------------------------------------
$objetivosPt = @'
A  disciplina  visa  apresentar  aos  estudantes  conceitos  e  técnicas  essenciais  para  o desenvolvimento e gestão da criatividade visando a solução de problemas e a promoção da inovação sistemática.
'@
$objectivesEn = @'
The course aims to introduce the students to the basic concepts and techniques to the creativity development and management aiming the problems solving and the promotion of systematic innovation.
'@
$professorName = @'
5840820 - Gustavo Aristides Santana Martinez
'@
$programaResumidoPt = @'
Solução de problema. Fundamentos de TRIZ. Contradições técnicas e físicas. Análise de recursos. Padrões  inventivos  e  modelo  substancia-campo.  Algoritmo  de  resolução  de  problemas inventivos. Evolução de sistemas. Ferramentas para vencer a inercia mental. Processo de resolução de problemas.
'@
$shortSyllabusEn = @'
Problem solving. TRIZ fundamentals. Technical and physical contradictions. Resources analysis. Inventive patterns and substance field model. Resolution algorithm of inventive problems. Systems evolution. Tools to overcome the mental inertia. Process of problem solving.
'@
$programaPt = @'
1) Enfoques tradicionais de solução de problema: tentativa e erro, brainstorming, seis chapéus.
2) Fundamentos de TRIZ: histórico, bases teóricas, níveis de inovação, ferramentas, estruturação sistemática, idealidade do sistema.
3) Contradições técnicas e físicas: conceito de contradição, contradição física, contradição técnica, matriz de contradições.
4) Análise de recursos: substancia, campo energético, espaço, tempo, informação, recursos funcionais.
5) Padrões inventivos e modelo substancia-campo: regras, solução padrão.
6) Algoritmo de resolução de problemas inventivos (ARIZ): modelo de problemas, passos.
7) Evolução de sistemas: leis de evolução, leis estáticas, leis cinemáticas, leis dinâmicas.
8) Ferramentas para vencer a inercia mental: 9 janelas, pequenos homens, sistemas alternativos, operadores de tamanho, tempo e custo.
9) Processo de resolução de problemas.
'@
$syllabusEn = @'
1) Traditional problem solving approaches: trial and error, brainstorming, Six Hats
2) TRIZ fundamentals: history, theoretical bases, levels of innovation, tools, systematic structure, and system ideality.
3) Technical and physical contradictions: concept of contradiction, physical contradiction, technical contradiction, contradictions matrix.
4) Resources Analysis: substance, energy field, space, time, information, functional resources.
5) Inventive patterns and substance field model: rules, standard solution.
6) Inventive problems solving algorithm (ARIZ): problems model, steps.
7) Systems evolution: laws of evolution, static laws, kinematic laws, dynamic laws.
8) Tools to overcome the mental inertia: nine windows, little men, alternative systems, size operators, time and cost.
9) Problems solving process.
'@
$metodoTexto = @'
Aulas expositivas, discussão de casos em sala de aula, painéis, debates, seminários, análise de vídeos e palestrantes externos.
'@
$criterioTexto = @'
NF=(N1 + N2)/2 Onde:
NF = nota final. 
N = nota
'@
$normaRecuperacaoTexto = @'
Estará apto a efetuar a prova de reavaliação o aluno que tiver como média final na disciplina uma nota igual ou superior a três (3,0) e inferior a cinco (5,0), e tiver, no mínimo, 70% (setenta por cento) de frequência às aulas. O cálculo de uma média aritmética simples será feito com a nota da prova de reavaliação e a média final obtida pelo aluno na disciplina. Se esta média resultar em nota igual ou superior a cinco (5,0), o aluno será aprovado.
'@
$bibliografiaTexto = @'
AZNAR, G. Ideias - 100 Técnicas de Criatividade. São Paulo: Summus. 2011. 256 p.
BARLACH, Lisete. "A Criatividade Humana sob a Ótica do Empreendedorismo Inovador". Tese de Doutorado. Sp: Instituto de Psicologia da Universidade de São Paulo, 2009. GURGEL, M.F. Criatividade & inovação: uma proposta de gestão da criatividade para o desenvolvimento da inovação. Rio de Janeiro: UFRJ, 2006, 193 p. (Dissertação. Mestrado em engenharia de produção).
SIQUEIRA, J. Criatividade Aplicada: habilidades e técnicas essenciais para a criatividade,
inovação e solução de problemas (3ª Edição). Siqueira Assessoria Empresarial. 2012. 113p. Disponível em < http://criatividadeaplicada.com/>.
PREDEBON, J. Criatividade (7ª Edição). São Paulo: Atlas. 2010. 246p.
UE (Unión Européia). Guía de buenas prácticas en materia de creatividad empresarial. Madrid: SUDOE.125p. (Publicado en el marco del proyecto Crea-Business-Idea del Programa SUDOE Interreg IV B).
 Bibliografia Complementar:
AMORIM, M.C.S.; FREDERICO, R. Criatividade, inovação e controle nas organizações. Revista de Ciências Humanas, Florianópolis, v. 42, n. 1 e 2, p. 75-89, 2008.
Baille, C., Enhancing Creativity in Science & Engineering, LTSN Composites Workshop, London, May 16, 2002.
Campbell, Brian, If TRIZ is Such a Good Idea, Why Isn't Everyone Using It?, The TRIZ Journal, April 2002.
Comments On: Campbell, Brian, If TRIZ is Such a Good Idea, Why Isn't Everyone Using It?, The TRIZ Journal, April 2002.
Filmore, Paul, The Real World: TRIZ in Two Hours for Undergraduate and Masters Level Students!, Proceedings of TRIZCON2006, Milwaukee, Wisconsin, U.S.A, April 4-5, 2006. Filmore, Paul, Thomond, P., Why Reinvent the Wheel? The Efficacy of Systematic Problem Solving Method TRIZ and Its Value for Innovation in Engineering and Its Implications for Engineering Management, Hong Kong Institute of Value Management, 7th International Conference, June 2005.
KELLEY, T.; LITTMAN, J. As 10 faces da inovação: estratégias para turbinar a criatividade. Rio de Janeiro: Campus/Elsevier. 2007. 280 p.
Kim, D.H., The Link Between Organizational and Individual Learning, Sloan Management Review, Fall 1993, pp 37-50.
Kolb, D.A., Experiential Learning. Experience as the Source of Learning and Development, Englewood Cliffs, NJ: Prentice-Hall, 1984.
Lumsdaine, E., Lumsdaine, M., Creative Problem Solving: Thinking Skills for a Changing World, McGraw-Hill, 1995.
Mann, Darrell, The Space Between Generic and Specific Problem Solutions, The TRIZ Journal, June 2001.
Mann, Darrell, Hands-On Systematic Innovation, CREAX Press, 2002, ISBN 90-77071-02-4. Mistry, J., White, F., and Berardi, A., Skills at Masters' Level in Geography Higher Education: Teaching, Learning and Applying, Planet 16, The Higher Education Academy: GEES, July 2006. NIZO, R. Foco e criatividade: fazer mais com menos. São Paulo: Summus. 2009. 120 p. RICKARDS, T. Creativity and problem solving at work. Gower: Aldershot, 1997. Robinson, M., Problem Solving in Groups, 2nd Ed, Gower, United Kingdom, 1993. Senge, P., The Fifth Discipline: the Art and Practice of the Learning Organisation, Century Business, 1990.
Wu, Tzann-Dwo, The Study of Problem Solving by TRIZ and Taguchi Methodology in Automobile Muffler Design, The TRIZ Journal, March 2004.
'@

# ---------------------------------------------------------------------------
# "Objetivos:" (row 10) now carries the Portuguese objectives paragraph;
# "Objectives:" (row 11) keeps the same English paragraph it already had.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt
$ws.Range("B11").Value = $objectivesEn
$ws.Range("C11").Value = $objectivesEn

# ---------------------------------------------------------------------------
# A brand-new row is inserted right after "Docentes responsaveis:" (row 12).
# Everything that used to start at row 13 shifts down to row 14 onward.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# Row 13: professor name only (column A left blank). Column B/C formatting is
# copied from the row below so no new cell style gets created.
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A13").Value = $null
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = $professorName
$ws.Range("C13").Value = $professorName
$ws.Rows.Item(13).RowHeight = 15

# Row 14: "Programa resumido:" now has real content instead of the stray "Semestral" placeholder
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = $programaResumidoPt
$ws.Range("C14").Value = $programaResumidoPt
$ws.Rows.Item(14).RowHeight = 60

# Row 15: "Short syllabus:" -- label/content unchanged, just re-asserted after the shift
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = $shortSyllabusEn
$ws.Range("C15").Value = $shortSyllabusEn
$ws.Rows.Item(15).RowHeight = 60

# Row 16: "Programa:" now has the full Portuguese syllabus instead of being empty
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt
$ws.Rows.Item(16).RowHeight = 120

# Row 17: "Syllabus:" -- label/content unchanged, just re-asserted after the shift
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = $syllabusEn
$ws.Range("C17").Value = $syllabusEn
$ws.Rows.Item(17).RowHeight = 120

# Row 18: "Avaliação:" section header, label only
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Rows.Item(18).RowHeight = 15

# Row 19: "Método:" now has real content instead of being empty
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = $metodoTexto
$ws.Range("C19").Value = $metodoTexto
$ws.Rows.Item(19).RowHeight = 60

# Row 20: "Critério:" with the NF=(N1+N2)/2 grading formula description
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = $criterioTexto
$ws.Range("C20").Value = $criterioTexto
$ws.Rows.Item(20).RowHeight = 60

# Row 21: "Norma de recuperação:" now has its own explanatory paragraph
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = $normaRecuperacaoTexto
$ws.Range("C21").Value = $normaRecuperacaoTexto
$ws.Rows.Item(21).RowHeight = 60

# Row 22 (brand new): "Bibliografia:" with the full reading list
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = $bibliografiaTexto
$ws.Range("C22").Value = $bibliografiaTexto
$ws.Rows.Item(22).RowHeight = 120
